$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be forced to Text format first,
# otherwise Excel auto-converts the assigned string into a number and the
# original text formatting (e.g. "531.72") would be lost. We reset the style
# back to Normal immediately after so no stray formatting is left behind.

$ws.Range("D2").Value = "57.885.37"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.274.78"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.17%  "
$ws.Range("D9").Value = "2.274.46"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0992"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.330"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "2.681.53"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "57.854.64"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "2.283.41"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  -1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("D31").Value = "0.0₃0716"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0944"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.550"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0209"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
